$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (pushes "Sorted Square" and everything below
# down by one row). The new row inherits formatting from the row above,
# matching Excel's default Insert behaviour.
$ws.Rows("7").Insert()

# --- Fill in the new row 7: "Running Sum of 1d Array" ---
$ws.Cells.Item(7,1).Value = "Arrays"
$ws.Cells.Item(7,2).Value = "Running Sum of 1d Array"
$ws.Cells.Item(7,3).Value = "https://leetcode.com/problems/running-sum-of-1d-array/"
$ws.Cells.Item(7,4).Value = "Using Separate Array, Using Input Array for Output"

$note = "Their's two approach to solve this challenge , either create a seprate array to store sum of elements (result[i] = result[i - 1] + nums[i];) , or use same input array for return output (nums[i] += nums[i - 1];)"
$ws.Cells.Item(7,5).Value = $note

$e7 = $ws.Cells.Item(7,5)

$run2 = $e7.Characters(104, 35)
$run2.Font.Name = "Arial (Body)"
$run2.Font.Size = 14
$run2.Font.Color = 255

$run3 = $e7.Characters(139, 46)
$run3.Font.Name = "Arial (Body)"
$run3.Font.Size = 14
$run3.Font.Color = 0

$run4 = $e7.Characters(185, 26)
$run4.Font.Name = "Arial (Body)"
$run4.Font.Size = 14
$run4.Font.Color = 255

# --- Trim the leading space from the two "Number of Ways to Split Array" names ---
$ws.Cells.Item(15,2).Value = "Number of Ways to Split Array (Without Array)"
$ws.Cells.Item(16,2).Value = "Number of Ways to Split Array (With Array)"

# --- Rebuild hyperlinks: the original hyperlink anchors shifted down by one
# row (from row 7 onward) because of the inserted row, so recreate the whole
# collection against the now-correct cells, then add the new one for row 7. ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C3"), "https://leetcode.com/problems/two-sum/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://leetcode.com/problems/valid-palindrome") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://leetcode.com/problems/merge-sorted-array") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://leetcode.com/problems/is-subsequence") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://leetcode.com/problems/reverse-string") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C8"), "https://leetcode.com/problems/squares-of-a-sorted-array/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C11"), "https://leetcode.com/problems/subarray-product-less-than-k/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C13"), "https://leetcode.com/problems/maximum-average-subarray-i/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C14"), "https://leetcode.com/problems/max-consecutive-ones-iii/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C15"), "https://leetcode.com/problems/number-of-ways-to-split-array/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C16"), "https://leetcode.com/problems/number-of-ways-to-split-array/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://leetcode.com/problems/running-sum-of-1d-array/") | Out-Null

# --- Update the active selection to match the saved view state ---
$ws.Range("B2").Select()
